# First draft layout of HV circuit: strip the old samtec-adapter-v2a BOM
# rows (components, values, vendors, part numbers, and their hyperlinks)
# out of the sheet, leaving just the header / totals scaffolding and the
# one remaining "Samtec QRM078" component note, then move the active
# selection onto the now-empty component-type column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all of the old BOM line-item data (values, component labels,
# types, vendors, part numbers, and the "link"/Digikey hyperlink column)
# for rows 3 through 10, but keep the M/O column formatting (Hyperlink
# style / date style) intact on the now-empty cells.
$ws.Range("D3:M10").ClearContents()

# Drop the now-orphaned hyperlink objects (Digikey / RFMW product links).
$ws.Hyperlinks.Delete()

# The only remaining populated note in that block switches from the RF
# switch part to the Samtec connector.
$ws.Range("I2").Value = "Samtec QRM078"

# Move the on-screen selection to reflect where editing continued.
$ws.Range("I3").Select()
